$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix genus name typo: Holzapfelia -> Holzapfeliella
$ws.Range("A3").Value = "Holzapfeliella"

# Insert a new row for genus "Paralactobacillus" above "Latilactobacillus" (row 12)
$ws.Rows(12).Insert()
$ws.Range("A12").Value = "Paralactobacillus"
$ws.Range("B12").Value = "closely related to lactobacilli"
$ws.Range("C12").Value = "P. selangorensis"
$ws.Range("D12").Value = "Homofermentative, vancomycin resistant, mesophilic organism."

# Fix typo in Limosilactobacillus type species (now shifted down to row 20)
$ws.Range("C20").Value = "L  fermentum"

# Insert two new rows for genera "Philodulcilactobacillus" and "Nicoliella"
# above "Apilactobacillus" (now shifted down to row 23)
$ws.Rows("23:24").Insert()

$ws.Range("A23").Value = "Philodulcilactobacillus"
$ws.Range("B23").Value = "Sugar-loving lactobacilli"
$ws.Range("C23").Value = "Pl. myokoensis"
$ws.Range("D23").Value = "Heterofermentative, vancomycin resistant, fructiphlic, growth on media solidified with gellan gum but not with agar"

$ws.Range("A24").Value = "Nicoliella"
$ws.Range("B24").Value = "Nicola Spurrier’s lactobacilli"
$ws.Range("C24").Value = "N. spurrieriana"
$ws.Range("D24").Value = "Heterofermentative, vancomycin resistant, small genome size but the single Nicoliella species has a broader carbohydrate fermentation pattern and a larger genome size than the closely related apilactobacilli. Likely adapted to bees or flowers."
$ws.Range("E24").Value = "10.1099/ijsem.0.005588"
